
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: Neg_Change
# ---------------------------------------------------------------
$wsNeg = $wb.Worksheets.Item("Neg_Change")

$negData = @(
    @("GRASIM", 2738, 2789, 2728.7, 2762, 489522, 1101889, -0.5557429105835524, "GRASIM"),
    @("HINDALCO", 792, 799.9, 783.85, 784.9, 5914226, 13189011, -0.5515792654961013, "HINDALCO"),
    @("DIVISLAB", 6700, 6748, 6636.5, 6706, 418569, 855609, -0.5107940659810731, "DIVISLAB"),
    @("ICICIGI", 2020.2, 2045, 2001, 2021.5, 431101, 882415, -0.5114532277896454, "ICICIGI"),
    @("ADANIENSOL", 965.55, 972.5, 955.1, 959, 674839, 1384632, -0.512622126312262, "ADANIENSOL"),
    @("RECLTD", 365.25, 367.15, 362, 362.75, 2951020, 6342222, -0.5347025064717066, "RECLTD"),
    @("PFC", 378, 381.65, 374.3, 377.4, 7558172, 16021457, -0.5282469003911442, "PFC"),
    @("LICI", 926.9, 931.9, 898, 899.9, 3204796, 7699680, -0.5837754296282442, "LICI"),
    @("JUBLFOOD", 577.95, 588.9, 575.5, 587, 866793, 1703407, -0.4911415768515687, "JUBLFOOD"),
    @("NHPC", 80.90000000000001, 82.72, 80.68000000000001, 82.52, 10235432, 21856711, -0.5317030087463754, "NHPC"),
    @("GODREJPROP", 2140, 2182, 2131.1, 2166.7, 276221, 591444, -0.532971845178918, "GODREJPROP"),
    @("HDFCAMC", 5432.5, 5510, 5416.5, 5490, 184862, 454966, -0.5936795277009711, "HDFCAMC"),
    @("POLICYBZR", 1783.8, 1812.9, 1761.8, 1800, 800584, 1906537, -0.5800847295384249, "POLICYBZR"),
    @("LUPIN", 1980.1, 2028, 1980.1, 1985, 1717376, 3553381, -0.5166924120999127, "LUPIN"),
    @("HINDPETRO", 479.55, 485.2, 477.75, 481.6, 2222672, 5059893, -0.5607274699287119, "HINDPETRO"),
    @("UPL", 750.05, 759, 746.7, 750.1, 2547152, 5598140, -0.5450003036722911, "UPL"),
    @("IDEA", 9.619999999999999, 9.74, 9.48, 9.529999999999999, 536169122, 1075464880, -0.5014536206891294, "IDEA"),
    @("PAYTM", 1340, 1353.8, 1325.4, 1330, 2677468, 5871932, -0.5440226487636437, "PAYTM"),
    @("IEX", 139.65, 140.8, 138.85, 139.5, 3498788, 8004546, -0.5628998821419728, "IEX"),
    @("DELHIVERY", 429.95, 430.75, 423.05, 428.1, 3333571, 6553523, -0.4913314563785006, "DELHIVERY"),
    @("BANDHANBNK", 154.39, 155.27, 152.25, 152.38, 7260687, 16575533, -0.561963588139217, "BANDHANBNK"),
)

# Clear out the old data rows (below header) before writing the new set,
# since the new data has a different number of rows than before.
$oldUsedRange = $wsNeg.UsedRange
$oldLastRow = $oldUsedRange.Rows.Count
if ($oldLastRow -gt 1) {
    $wsNeg.Range($wsNeg.Cells.Item(2, 1), $wsNeg.Cells.Item($oldLastRow, 9)).Clear()
}

for ($i = 0; $i -lt $negData.Count; $i++) {
    $row = $negData[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsNeg.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

$wsNeg.Range("A1:I22").EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------
# Sheet 2: Pos_Change
# ---------------------------------------------------------------
$wsPos = $wb.Worksheets.Item("Pos_Change")

$posData = @(
    @("WIPRO", 237.44, 241.55, 237.02, 239.85, 8566758, 5520619, 0.5517749006044431, "WIPRO"),
    @("KOTAKBANK", 2081, 2104.8, 2081, 2090, 2173296, 1400189, 0.5521447461735522, "KOTAKBANK"),
    @("SIEMENS", 3033.4, 3088, 3003.1, 3053, 465289, 304289, 0.5291022679097831, "SIEMENS"),
    @("VOLTAS", 1320, 1374.8, 1308.5, 1340.5, 1057635, 666465, 0.5869325470954964, "VOLTAS"),
    @("MPHASIS", 2736, 2804, 2726.6, 2775, 430038, 268774, 0.5999985117608102, "MPHASIS"),
    @("HUDCO", 231, 236.45, 230.2, 233.32, 4276399, 3027164, 0.4126750318119534, "HUDCO"),
    @("BHARATFORG", 1318.3, 1336, 1306.8, 1330.2, 723187, 457443, 0.5809335807958588, "BHARATFORG"),
    @("ASHOKLEY", 141.8, 143.15, 140.69, 142.41, 7609145, 4888203, 0.5566344114595896, "ASHOKLEY"),
    @("SJVN", 85.27, 85.40000000000001, 83.23999999999999, 84.18000000000001, 3617159, 2376356, 0.5221452509640812, "SJVN"),
    @("LAURUSLABS", 985, 1004, 983, 1000.1, 1569187, 1087182, 0.4433526309302399, "LAURUSLABS"),
)

$oldUsedRangePos = $wsPos.UsedRange
$oldLastRowPos = $oldUsedRangePos.Rows.Count
if ($oldLastRowPos -gt 1) {
    $wsPos.Range($wsPos.Cells.Item(2, 1), $wsPos.Cells.Item($oldLastRowPos, 9)).Clear()
}

for ($i = 0; $i -lt $posData.Count; $i++) {
    $row = $posData[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsPos.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

$wsPos.Range("A1:I11").EntireColumn.AutoFit() | Out-Null

